$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Fill in row 61 (previously empty) with a new work session entry
$ws.Range("C61").Value = 0.33333333333333331
$ws.Range("D61").Value = "-"
$ws.Range("E61").Value = 0.39583333333333331
$ws.Range("F61").Value = "Rédaction rapport"

# Fill in row 62 (previously empty) with a new work session entry
$ws.Range("C62").Value = 0.39583333333333331
$ws.Range("D62").Value = "-"
$ws.Range("E62").Value = 0.4375
$ws.Range("F62").Value = "Modularisation du code, réflexions sur la partie Archers"

# Force a full recalculation so dependent formulas (E66, H5) update their cached values
$excel.CalculateFullRebuild()

# Update the active selection on the sheet to match the edited cell
$ws.Range("C63").Select()
